$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 832.05884
$ws.Range("I6").Value = 280.625
$ws.Range("J6").Value = 1322.2222
$ws.Range("K6").Value = 841.875
$ws.Range("L6").Value = 3966.6666
$ws.Range("M6").Value = -729.875
$ws.Range("N6").Value = -4190.6666

$ws.Range("H9").Value = 90
$ws.Range("I9").Value = 53.333332
$ws.Range("K9").Value = 53.333332
$ws.Range("M9").Value = 115.666668

$ws.Range("H41").Value = 731.9091
$ws.Range("I41").Value = 333
$ws.Range("J41").Value = 771.8
$ws.Range("K41").Value = 333
$ws.Range("L41").Value = 771.8
$ws.Range("M41").Value = 107
$ws.Range("N41").Value = -1651.8

$ws.Range("H53").Value = 53293.684
$ws.Range("I53").Value = 101075.1
$ws.Range("J53").Value = 203.22223
$ws.Range("K53").Value = 101075.1
$ws.Range("L53").Value = 203.22223
$ws.Range("M53").Value = -100438.1
$ws.Range("N53").Value = -1477.22223

$ws.Range("H62").Value = 1163.7273
$ws.Range("I62").Value = 1275.625
$ws.Range("J62").Value = 865.3333
$ws.Range("K62").Value = 1275.625
$ws.Range("L62").Value = 865.3333
$ws.Range("M62").Value = -651.625
$ws.Range("N62").Value = -2113.3333

$ws.Range("H65").Value = 1163.7273
$ws.Range("I65").Value = 1275.625
$ws.Range("J65").Value = 865.3333
$ws.Range("K65").Value = 6378.125
$ws.Range("L65").Value = 4326.6665
$ws.Range("M65").Value = -3258.125
$ws.Range("N65").Value = -10566.6665

$ws.Range("H138").Value = 6619.9644
$ws.Range("I138").Value = 1156.7273
$ws.Range("J138").Value = 7955.4224
$ws.Range("K138").Value = 3470.1819
$ws.Range("L138").Value = 23866.2672
$ws.Range("M138").Value = 1669.8181
$ws.Range("N138").Value = -34146.2672

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 7000
$ws.Range("I3").Value = 7000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 7000
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -6885
$ws.Range("N3").ClearContents()

$ws.Range("H45").Value = 9116.571
$ws.Range("I45").Value = 10293.909
$ws.Range("J45").Value = 4799.6665
$ws.Range("K45").Value = 10293.909
$ws.Range("L45").Value = 4799.6665
$ws.Range("M45").Value = -9916.909
$ws.Range("N45").Value = -5553.6665

$ws.Range("H122").Value = 1222492.2
$ws.Range("I122").Value = 1350938.9
$ws.Range("J122").Value = 2250
$ws.Range("K122").Value = 4052816.7
$ws.Range("L122").Value = 6750
$ws.Range("M122").Value = -4050366.7
$ws.Range("N122").Value = -11650

$ws.Range("H123").Value = 29999
$ws.Range("J123").Value = 29999
$ws.Range("L123").Value = 29999
$ws.Range("N123").Value = -39799

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 197.5
$ws.Range("I12").Value = 96.666664
$ws.Range("J12").Value = 500
$ws.Range("K12").Value = 96.666664
$ws.Range("L12").Value = 500
$ws.Range("M12").Value = 71.333336
$ws.Range("N12").Value = -836

$ws.Range("H86").Value = 2213.25
$ws.Range("I86").Value = 2201
$ws.Range("J86").Value = 2250
$ws.Range("K86").Value = 2201
$ws.Range("L86").Value = 2250
$ws.Range("M86").Value = -1078
$ws.Range("N86").Value = -4496

$ws.Range("H89").Value = 2213.25
$ws.Range("I89").Value = 2201
$ws.Range("J89").Value = 2250
$ws.Range("K89").Value = 11005
$ws.Range("L89").Value = 11250
$ws.Range("M89").Value = -5389
$ws.Range("N89").Value = -22482

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 11618.28
$ws.Range("J60").Value = 11618.28
$ws.Range("L60").Value = 11618.28
$ws.Range("N60").Value = -12640.28

$ws.Range("H62").Value = 12504
$ws.Range("I62").Value = 10003.333
$ws.Range("K62").Value = 10003.333
$ws.Range("M62").Value = -9379.333000000001

$ws.Range("H65").Value = 12504
$ws.Range("I65").Value = 10003.333
$ws.Range("K65").Value = 50016.665
$ws.Range("M65").Value = -46896.665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 687900
$ws.Range("I4").Value = 786100
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 2358300
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = -2358188
$ws.Range("N4").Value = -1724

$ws.Range("H9").Value = 133334060
$ws.Range("I9").Value = 333333630
$ws.Range("J9").Value = 83334180
$ws.Range("K9").Value = 1000000890
$ws.Range("L9").Value = 250002540
$ws.Range("M9").Value = -1000000666
$ws.Range("N9").Value = -250002988

$ws.Range("H16").Value = 140
$ws.Range("I16").Value = 140
$ws.Range("K16").Value = 420
$ws.Range("M16").Value = -247

$ws.Range("H17").Value = 1335.7142
$ws.Range("I17").Value = 1220
$ws.Range("J17").Value = 1625
$ws.Range("K17").Value = 3660
$ws.Range("L17").Value = 4875
$ws.Range("M17").Value = -3491
$ws.Range("N17").Value = -5213

$ws.Range("H19").Value = 3000
$ws.Range("J19").Value = 3000
$ws.Range("L19").Value = 9000
$ws.Range("N19").Value = -9348

$ws.Range("H25").Value = 3100.1667
$ws.Range("I25").Value = 1766.3334
$ws.Range("J25").Value = 4434
$ws.Range("K25").Value = 5299.0002
$ws.Range("L25").Value = 13302
$ws.Range("M25").Value = -5130.0002
$ws.Range("N25").Value = -13640

$ws.Range("H30").Value = 3100.1667
$ws.Range("I30").Value = 1766.3334
$ws.Range("J30").Value = 4434
$ws.Range("K30").Value = 5299.0002
$ws.Range("L30").Value = 13302
$ws.Range("M30").Value = -5197.0002
$ws.Range("N30").Value = -13506

$ws.Range("H122").Value = 1004.1429
$ws.Range("I122").Value = 866.3333
$ws.Range("J122").Value = 1107.5
$ws.Range("K122").Value = 7796.9997
$ws.Range("L122").Value = 9967.5
$ws.Range("M122").Value = -5346.9997
$ws.Range("N122").Value = -14867.5

$ws.Range("H131").Value = 23415754
$ws.Range("I131").Value = 7731197
$ws.Range("J131").Value = 29412790
$ws.Range("K131").Value = 23193591
$ws.Range("L131").Value = 88238370
$ws.Range("M131").Value = -23188551
$ws.Range("N131").Value = -88248450

$ws.Range("H133").Value = 28855.412
$ws.Range("J133").Value = 8941.379000000001
$ws.Range("L133").Value = 26824.137
$ws.Range("N133").Value = -36944.137

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 40000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 40000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H122").Value = 1908458.9
$ws.Range("I122").Value = 3088366.5
$ws.Range("J122").Value = 2454.3076
$ws.Range("K122").Value = 9265099.5
$ws.Range("L122").Value = 7362.9228
$ws.Range("M122").Value = -9262649.5
$ws.Range("N122").Value = -12262.9228

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 18518878
$ws.Range("I55").Value = 271.1111
$ws.Range("J55").Value = 27778182
$ws.Range("K55").Value = 271.1111
$ws.Range("L55").Value = 27778182
$ws.Range("M55").Value = -98.11110000000002
$ws.Range("N55").Value = -27778528

$ws.Range("H60").Value = 20000
$ws.Range("I60").Value = 10000
$ws.Range("J60").Value = 30000
$ws.Range("K60").Value = 10000
$ws.Range("L60").Value = 30000
$ws.Range("M60").Value = -9491
$ws.Range("N60").Value = -31018

$ws.Range("H82").Value = 1838950
$ws.Range("I82").Value = 3334233.2
$ws.Range("J82").Value = 343666.66
$ws.Range("K82").Value = 3334233.2
$ws.Range("L82").Value = 343666.66
$ws.Range("M82").Value = -3333872.2
$ws.Range("N82").Value = -344388.66

$ws.Range("H85").Value = 1838950
$ws.Range("I85").Value = 3334233.2
$ws.Range("J85").Value = 343666.66
$ws.Range("K85").Value = 3334233.2
$ws.Range("L85").Value = 343666.66
$ws.Range("M85").Value = -3332985.2
$ws.Range("N85").Value = -346162.66

$ws.Range("H122").Value = 3704954
$ws.Range("I122").Value = 4766059
$ws.Range("J122").Value = 1431157.1
$ws.Range("K122").Value = 14298177
$ws.Range("L122").Value = 4293471.300000001
$ws.Range("M122").Value = -14295727
$ws.Range("N122").Value = -4298371.300000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 14001
$ws.Range("I62").Value = 4334.3335
$ws.Range("J62").Value = 23667.666
$ws.Range("K62").Value = 4334.3335
$ws.Range("L62").Value = 23667.666
$ws.Range("M62").Value = -3710.3335
$ws.Range("N62").Value = -24915.666

$ws.Range("H65").Value = 14001
$ws.Range("I65").Value = 4334.3335
$ws.Range("J65").Value = 23667.666
$ws.Range("K65").Value = 21671.6675
$ws.Range("L65").Value = 118338.33
$ws.Range("M65").Value = -18551.6675
$ws.Range("N65").Value = -124578.33

$ws.Range("H81").Value = 1325.5
$ws.Range("I81").Value = 900
$ws.Range("J81").Value = 1751
$ws.Range("K81").Value = 1800
$ws.Range("L81").Value = 3502
$ws.Range("M81").Value = -739
$ws.Range("N81").Value = -5624

$ws.Range("H84").Value = 1325.5
$ws.Range("I84").Value = 900
$ws.Range("J84").Value = 1751
$ws.Range("K84").Value = 9000
$ws.Range("L84").Value = 17510
$ws.Range("M84").Value = -3696
$ws.Range("N84").Value = -28118

$ws.Range("H122").Value = 2164
$ws.Range("I122").Value = 1779.5
$ws.Range("J122").Value = 2933
$ws.Range("K122").Value = 5338.5
$ws.Range("L122").Value = 8799
$ws.Range("M122").Value = -2888.5
$ws.Range("N122").Value = -13699

$ws.Range("H123").Value = 29749.857
$ws.Range("J123").Value = 29749.857
$ws.Range("L123").Value = 29749.857
$ws.Range("N123").Value = -39549.857
